# Regenerate merged AHB files
#
# 1. Header row: rename the "_old" suffixed column headers (A1:J1) to
#    "_FV2410" and the "_new" suffixed column headers (L1:U1) to
#    "_FV2504". Column K1 ("diff") is left untouched.
# 2. Turn the A1:U55 range into a real Excel Table ("Table1") with an
#    AutoFilter.
# 3. Freeze the header row (split/freeze below row 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns 1-10 (A-J) carry the "..._old" headers -> rename to "..._FV2410".
for ($col = 1; $col -le 10; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = ($cell.Value2.ToString() -replace "_old$", "_FV2410")
}

# Column 11 (K) is "diff" and is left untouched.

# Columns 12-21 (L-U) carry the "..._new" headers -> rename to "..._FV2504".
for ($col = 12; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = ($cell.Value2.ToString() -replace "_new$", "_FV2504")
}

# Convert the used range into an Excel Table with an AutoFilter.
$tableRange = $ws.Range("A1:U55")
$tbl = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$tbl.Name = "Table1"

# Freeze the header row.
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
